# Weekly update: insert a new price-report row for
# "Feria Lagunitas de Puerto Montt - Membrillo" (week of 2023-03-28,
# Excel serial date 45013), pushing the existing rows 128-131 down to
# 129-132.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 128; this shifts the previous
# rows 128..131 down to 129..132 and expands the used range / dimension.
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with this week's data.
$ws.Cells.Item(128, 1).Value  = 4
$ws.Cells.Item(128, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(128, 3).Value  = "Los Lagos"
$ws.Cells.Item(128, 4).Value  = 45013
$ws.Cells.Item(128, 5).Value  = 10
$ws.Cells.Item(128, 6).Value  = "Fruta"
$ws.Cells.Item(128, 7).Value  = 100104
$ws.Cells.Item(128, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(128, 9).Value  = 100104003
$ws.Cells.Item(128, 10).Value = "Membrillo"
$ws.Cells.Item(128, 11).Value = "Champion"
$ws.Cells.Item(128, 12).Value = "Primera"
$ws.Cells.Item(128, 13).Value = 300
$ws.Cells.Item(128, 14).Value = 17000
$ws.Cells.Item(128, 15).Value = 18000
$ws.Cells.Item(128, 16).Value = 17500
$ws.Cells.Item(128, 17).Value = "$/caja 18 kilos empedrada"
$ws.Cells.Item(128, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(128, 19).Value = 972
$ws.Cells.Item(128, 20).Value = 18
